# Update the "取得日時" (retrieved timestamp) column for all data rows
# on the first worksheet to reflect the latest scrape time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2025-09-18 06:32:45"

for ($row = 2; $row -le 10; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
